$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 127
$ws1.Range("F4").Value = 896
$ws1.Range("F5").Value = 1068
$ws1.Range("F6").Value = 1547
$ws1.Range("F8").Value = 668
$ws1.Range("F9").Value = 12095
$ws1.Range("F11").Value = 2154
$ws1.Range("F13").Value = 251
$ws1.Range("F15").Value = 1221
$ws1.Range("F16").Value = 199
$ws1.Range("F17").Value = 268
$ws1.Range("F18").Value = 770
$ws1.Range("F19").Value = 672
$ws1.Range("F20").Value = 292
$ws1.Range("F21").Value = 2922
$ws1.Range("F22").Value = 753
$ws1.Range("F23").Value = 4022
$ws1.Range("F24").Value = 1108
$ws1.Range("F25").Value = 856
$ws1.Range("F29").Value = 1038
$ws1.Range("F30").Value = 47
$ws1.Range("F31").Value = 95
$ws1.Range("F32").Value = 270
$ws1.Range("F34").Value = 19
$ws1.Range("F36").Value = 10
$ws1.Range("F37").Value = 4403
$ws1.Range("F38").Value = 14
$ws1.Range("F39").Value = 4524
$ws1.Range("F40").Value = 5529
$ws1.Range("F42").Value = 128
$ws1.Range("F43").Value = 60
$ws1.Range("F44").Value = 170
$ws1.Range("F45").Value = 311
$ws1.Range("F47").Value = 40
$ws1.Range("F48").Value = 4109
$ws1.Range("F49").Value = 122

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 4171
$ws2.Range("F5").Value = 59
$ws2.Range("F6").Value = 98
$ws2.Range("F11").Value = 109
$ws2.Range("F13").Value = 1025
$ws2.Range("F19").Value = 47

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 760
$ws3.Range("F3").Value = 446
$ws3.Range("E3").Value = "2024.09.15 00:00-10.31 23:59"
$ws3.Range("I3").Value = "//i1.hdslb.com/bfs/openplatform/202409/RVDH3aey1727420551197.jpeg"
$ws3.Range("F4").Value = 73

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 760
$ws4.Range("F3").Value = 446
$ws4.Range("E3").Value = "2024.09.15 00:00-10.31 23:59"
$ws4.Range("I3").Value = "//i1.hdslb.com/bfs/openplatform/202409/RVDH3aey1727420551197.jpeg"
$ws4.Range("F5").Value = 73
$ws4.Range("F6").Value = 127
$ws4.Range("F8").Value = 896
$ws4.Range("F9").Value = 1068
$ws4.Range("F10").Value = 1547
$ws4.Range("F12").Value = 668
$ws4.Range("F13").Value = 12095
$ws4.Range("F14").Value = 2154
$ws4.Range("F16").Value = 1221
$ws4.Range("F17").Value = 199
$ws4.Range("F18").Value = 268
$ws4.Range("F19").Value = 770
$ws4.Range("F20").Value = 672
$ws4.Range("F21").Value = 2922
$ws4.Range("F22").Value = 753
$ws4.Range("F23").Value = 4023
$ws4.Range("F24").Value = 4023
$ws4.Range("F25").Value = 1108
$ws4.Range("F26").Value = 856
$ws4.Range("F32").Value = 1038
$ws4.Range("F33").Value = 47
$ws4.Range("F34").Value = 95
$ws4.Range("F35").Value = 109
$ws4.Range("F36").Value = 270
$ws4.Range("F38").Value = 4403
$ws4.Range("F41").Value = 128
$ws4.Range("F42").Value = 170
$ws4.Range("F43").Value = 311
$ws4.Range("F47").Value = 4109
$ws4.Range("F48").Value = 47
